$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# Row 10: new "MANY TO MANY" heading next to existing content
$ws.Range("G10").Value = "MANY TO MANY"

# Row 12: table headers for the many-to-many example (left to right)
$ws.Range("G12").Value = "USERS"
$ws.Range("J12").Value = "CLAIMS"
$ws.Range("M12").Value = "USERCLAIMS"

# USERS table: header row then data, filled column by column
$ws.Range("G13").Value = "NAME"
$ws.Range("G14").Value = "MICHAEL"
$ws.Range("G15").Value = "GABRIEL"
$ws.Range("G16").Value = "AMANDA"

# CLAIMS table: header row then data, filled column by column
$ws.Range("J13").Value = "TITLE"
$ws.Range("J14").Value = "ACCOUNT ACC"
$ws.Range("J15").Value = "SALES ACC"
$ws.Range("J16").Value = "PURCHASE ACC"

# USERCLAIMS table headers
$ws.Range("M13").Value = "USERID"
$ws.Range("N13").Value = "CLAIMSID"

# ID columns (F, I, L) reuse the existing "ID" shared string, right aligned
$ws.Range("F13").Value = "ID"
$ws.Range("F13").HorizontalAlignment = -4152
$ws.Range("I13").Value = "ID"
$ws.Range("I13").HorizontalAlignment = -4152
$ws.Range("L13").Value = "ID"
$ws.Range("L13").HorizontalAlignment = -4152

# Numeric ID columns for each table
$ws.Range("F14").Value = 1
$ws.Range("F15").Value = 2
$ws.Range("F16").Value = 3

$ws.Range("I14").Value = 1
$ws.Range("I15").Value = 2
$ws.Range("I16").Value = 3

$ws.Range("L14").Value = 1
$ws.Range("L15").Value = 2
$ws.Range("L16").Value = 3
$ws.Range("L17").Value = 4
$ws.Range("L18").Value = 5
$ws.Range("L19").Value = 6
$ws.Range("L20").Value = 7

# USERCLAIMS join data (USERID / CLAIMSID)
$ws.Range("M14").Value = 1
$ws.Range("N14").Value = 1

$ws.Range("M15").Value = 1
$ws.Range("N15").Value = 2

$ws.Range("M16").Value = 2
$ws.Range("N16").Value = 1

$ws.Range("M17").Value = 2
$ws.Range("N17").Value = 2

$ws.Range("M18").Value = 2
$ws.Range("N18").Value = 3

$ws.Range("M19").Value = 3
$ws.Range("N19").Value = 1

$ws.Range("M20").Value = 3
$ws.Range("N20").Value = 3

# Update the selection to match the author's final cursor position
$ws.Range("P15").Select()
